$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.307.47"
$ws.Cells.Item(2, 5).Value = "  +7.42%  "
$ws.Cells.Item(3, 4).Value = "3.534.77"
$ws.Cells.Item(3, 5).Value = "  +11.00%  "
$ws.Cells.Item(4, 5).Value = "  +0.15%  "
$ws.Cells.Item(5, 4).Value = "'191.40"
$ws.Cells.Item(5, 5).Value = "  +11.05%  "
$ws.Cells.Item(6, 4).Value = "'553.74"
$ws.Cells.Item(6, 5).Value = "  +5.16%  "
$ws.Cells.Item(7, 4).Value = "3.524.65"
$ws.Cells.Item(7, 5).Value = "  +10.91%  "
$ws.Cells.Item(8, 4).Value = "'0.609"
$ws.Cells.Item(8, 5).Value = "  +2.51%  "
$ws.Cells.Item(9, 5).Value = "  -0.04%  "
$ws.Cells.Item(10, 4).Value = "'0.635"
$ws.Cells.Item(10, 5).Value = "  +5.10%  "
$ws.Cells.Item(11, 5).Value = "  +15.23%  "
$ws.Cells.Item(12, 4).Value = "'55.24"
$ws.Cells.Item(12, 5).Value = "  +3.97%  "
$ws.Cells.Item(13, 5).Value = "  +7.56%  "
$ws.Cells.Item(14, 4).Value = "'9.38"
$ws.Cells.Item(14, 5).Value = "  +3.53%  "
$ws.Cells.Item(15, 4).Value = "4.085.32"
$ws.Cells.Item(15, 5).Value = "  +10.80%  "
$ws.Cells.Item(16, 4).Value = "3.530.15"
$ws.Cells.Item(16, 5).Value = "  +11.44%  "
$ws.Cells.Item(17, 5).Value = "  +3.60%  "
$ws.Cells.Item(18, 4).Value = "67.328.58"
$ws.Cells.Item(18, 5).Value = "  +7.90%  "
$ws.Cells.Item(19, 4).Value = "'18.24"
$ws.Cells.Item(19, 5).Value = "  +6.07%  "
$ws.Cells.Item(20, 5).Value = "  +8.67%  "
$ws.Cells.Item(21, 5).Value = "  +3.39%  "
$ws.Cells.Item(22, 4).Value = "'432.45"
$ws.Cells.Item(22, 5).Value = "  +18.85%  "
$ws.Cells.Item(23, 5).Value = "  +4.33%  "
$ws.Cells.Item(24, 4).Value = "'84.97"
$ws.Cells.Item(24, 5).Value = "  +4.83%  "
$ws.Cells.Item(25, 5).Value = "  +7.71%  "
$ws.Cells.Item(26, 4).Value = "'11.23"
$ws.Cells.Item(26, 5).Value = "  +0.30%  "
$ws.Cells.Item(27, 4).Value = "'2.92"
$ws.Cells.Item(27, 5).Value = "  +10.92%  "
$ws.Cells.Item(28, 2).Value = "LEO"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(28, 4).Value = "'6.15"
$ws.Cells.Item(28, 5).Value = "  +0.68%  "
$ws.Cells.Item(29, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(29, 4).Value = "'12.04"
$ws.Cells.Item(29, 5).Value = "  +6.21%  "
$ws.Cells.Item(30, 2).Value = "Filecoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(30, 4).Value = "'8.99"
$ws.Cells.Item(30, 5).Value = "  +9.90%  "
$ws.Cells.Item(31, 2).Value = "EthereumClassic"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(31, 4).Value = "'30.33"
$ws.Cells.Item(31, 5).Value = "  +7.10%  "
$ws.Cells.Item(32, 2).Value = "Bittensor"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(32, 4).Value = "'649.63"
$ws.Cells.Item(32, 5).Value = "  +1.60%  "
$ws.Cells.Item(33, 2).Value = "NEARProtocol"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(33, 4).Value = "'6.70"
$ws.Cells.Item(33, 5).Value = "  +3.88%  "
$ws.Cells.Item(34, 2).Value = "Cosmos"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(34, 4).Value = "'11.75"
$ws.Cells.Item(34, 5).Value = "  +3.86%  "
$ws.Cells.Item(35, 2).Value = "Hedera"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(35, 4).Value = "'0.111"
$ws.Cells.Item(35, 5).Value = "  +5.70%  "
$ws.Cells.Item(36, 2).Value = "OKB"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(36, 4).Value = "'59.55"
$ws.Cells.Item(36, 5).Value = "  +5.64%  "
$ws.Cells.Item(37, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(37, 4).Value = "'38.81"
$ws.Cells.Item(37, 5).Value = "  +5.05%  "
$ws.Cells.Item(38, 2).Value = "PEPE"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(38, 4).Value = "0.0₃0821"
$ws.Cells.Item(38, 5).Value = "  +16.11%  "
$ws.Cells.Item(39, 2).Value = "Dai"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(39, 4).Value = "'0.999"
$ws.Cells.Item(39, 5).Value = "  -0.13%  "
$ws.Cells.Item(40, 2).Value = "TheGraph"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(40, 4).Value = "'0.391"
$ws.Cells.Item(40, 5).Value = "  +4.42%  "
$ws.Cells.Item(41, 2).Value = "Kaspa"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(41, 4).Value = "'0.143"
$ws.Cells.Item(41, 5).Value = "  +15.50%  "
$ws.Cells.Item(42, 2).Value = "Stacks"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(42, 4).Value = "'3.34"
$ws.Cells.Item(42, 5).Value = "  +14.46%  "
$ws.Cells.Item(43, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(43, 4).Value = "'0.999"
$ws.Cells.Item(43, 5).Value = "  +0.27%  "
$ws.Cells.Item(44, 2).Value = "Maker"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(44, 4).Value = "3.028.55"
$ws.Cells.Item(44, 5).Value = "  +5.18%  "
$ws.Cells.Item(45, 2).Value = "Fetch.AI"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(45, 4).Value = "'2.66"
$ws.Cells.Item(45, 5).Value = "  +5.65%  "
$ws.Cells.Item(46, 4).Value = "'2.90"
$ws.Cells.Item(46, 5).Value = "  +9.60%  "
$ws.Cells.Item(47, 2).Value = "ThetaToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(47, 4).Value = "'2.88"
$ws.Cells.Item(47, 5).Value = "  +11.55%  "
$ws.Cells.Item(48, 2).Value = "ApeXProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(48, 4).Value = "'3.35"
$ws.Cells.Item(48, 5).Value = "  +12.25%  "
$ws.Cells.Item(49, 2).Value = "VeChain"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(49, 4).Value = "'0.0419"
$ws.Cells.Item(49, 5).Value = "  +6.62%  "
$ws.Cells.Item(50, 2).Value = "Stellar"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(50, 4).Value = "'0.131"
$ws.Cells.Item(50, 5).Value = "  +6.20%  "
$ws.Cells.Item(51, 2).Value = "THORChain"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(51, 4).Value = "'8.76"
$ws.Cells.Item(51, 5).Value = "  +14.47%  "
